$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.338.38"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.786.06"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'32.91"
$ws.Range("E8").Value = "  +3.54%  "
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("D10").Value = "'0.0687"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").Value = "'0.0946"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "2.044.54"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "'11.18"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").Value = "1.795.63"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").Value = "'0.634"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "34.362.59"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("D18").Value = "'68.33"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "'244.91"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "0.0₃0793"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "'11.26"
$ws.Range("E21").Value = "  +3.44%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'4.15"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").Value = "'168.68"
$ws.Range("E24").Value = "  +4.43%  "
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("D26").Value = "'7.31"
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("D27").Value = "'16.51"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'4.02"
$ws.Range("E30").Value = "  +9.08%  "
$ws.Range("D31").Value = "'0.0525"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.79"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.23"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("D35").Value = "1.411.62"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("E36").Value = "  +5.15%  "
$ws.Range("D37").Value = "'0.681"
$ws.Range("E37").Value = "  +4.63%  "
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "'84.34"
$ws.Range("E40").Value = "  +5.08%  "
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.77"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.938"
$ws.Range("E43").Value = "  +2.52%  "
$ws.Range("D44").Value = "'14.01"
$ws.Range("E44").Value = "  +2.33%  "
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("D47").Value = "'6.07"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "1.945.73"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").Value = "'105.14"
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("E51").Value = "  -2.09%  "
